# Sync attendance_reports: reorder "Recorded By" names in column G
# so the last author listed moves to the front of the comma-separated list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$replacements = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
